# es_MX "Video Script_Take_a_Pause.docx" translation update
#
# Applies, in order, the same textual edits described by the commit's
# unified diff:
#   - "[pausa]"  -> "[pause]"          (every remaining occurrence)
#   - "En;"/"En; " -> "Inhalando;"/"Inhalando; "
#   - "y fuera;"/"y fuera; " -> "y exhalando;"/"y exhalando; "
#   - assorted sentence-level rewording in the "Take a Pause" table
#
# Each edit is scoped to the specific paragraph it belongs to (via
# $d.Paragraphs.Item(N).Range) so that duplicate strings elsewhere in the
# document are left untouched and run-level formatting is preserved.

$d = $word.ActiveDocument

$wdReplaceOne = 1

function Replace-InParagraph($index, $oldText, $newText) {
    $rng = $d.Paragraphs.Item($index).Range
    $ok = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceOne)
    if (-not $ok) {
        throw "Replacement failed in paragraph $index : '$oldText' -> '$newText'"
    }
}

Replace-InParagraph 16 "[pausa]" "[pause]"

Replace-InParagraph 19 "En;" "Inhalando;"
Replace-InParagraph 20 "[pausa]" "[pause]"
Replace-InParagraph 21 "y fuera; " "y exhalando; "

Replace-InParagraph 23 "En; " "Inhalando; "
Replace-InParagraph 24 "[pausa] " "[pause] "
Replace-InParagraph 25 "y fuera;" "y exhalando;"
Replace-InParagraph 26 "[pausa]" "[pause]"

Replace-InParagraph 27 "En; " "Inhalando; "
Replace-InParagraph 28 "[pausa] " "[pause] "
Replace-InParagraph 29 "y fuera;" "y exhalando;"
Replace-InParagraph 30 "[pausa]" "[pause]"

Replace-InParagraph 31 "Observa cómo se siente tu cuerpo mientras respiras. " "Siente tu cuerpo mientras respiras. "
Replace-InParagraph 32 "[pausa]" "[pause]"
Replace-InParagraph 33 "Observa dónde sientes tensión en el cuerpo." "Pon atención a dónde sientes tensión en el cuerpo."
Replace-InParagraph 34 "[pausa]" "[pause]"
Replace-InParagraph 35 "Intenta que tu cuerpo se relaje." "Deja que tu cuerpo se relaje."
Replace-InParagraph 36 "[pausa]" "[pause]"
Replace-InParagraph 37 "Cuando estés preparado, vuelve a abrir los ojos. " "Vuelve a abrir los ojos cuando estés listo. "

Replace-InParagraph 40 "Ahora, fíjate si te sientes diferente de cuando empezaste esta actividad." "Ahora, fíjate si te sientes diferente a cuando empezaste esta actividad."
Replace-InParagraph 41 "[pausa]" "[pause]"

Replace-InParagraph 44 "Intenta hacer una pausa cada vez que te sientas enfadado, abrumado, estresado o preocupado. " "Intenta Tomar una Pausa cada vez que te sientas enojado, abrumado, estresado o preocupado. "
Replace-InParagraph 46 "Incluso unas cuantas respiraciones profundas o la conexión con el suelo pueden marcar la diferencia. " "Hacer unas cuantas respiraciones profundas o conectar con el suelo debajo de ti pueden marcar la diferencia. "
Replace-InParagraph 48 "También puedes hacer una pausa con tu niña, niño o adolescente." "También puedes Tomar una Pausa con tu niña, niño o adolescente."

Write-Host "All replacements applied."
